# TC23_CDS_Filter_PHSAccession-phs002517_FileType_LibStrat_LibSrc.xlsx
# "CDS Test cases Study Facet and Many to Many test cases."
#
# The underlying OOXML diff for this workbook is dominated by Excel/engine
# generated metadata (fileVersion build number, revisionPtr GUID, the
# x15ac:absPath of the machine that last saved it, bookView window
# geometry, and sub-pixel column width / default row height drift coming
# from a different Excel build's font-metrics). None of that is reachable
# through the Excel object model - it is stamped by the application on
# save and is not an intentional content edit.
#
# The one genuinely user-driven change captured in the diff is the sheet
# view: the author scrolled the "startup" sheet back to the top
# (topLeftCell B4 -> B1) and left the selection on D3 (was C4) before
# saving.
#
# (The sharedStrings.xml churn - si index 12/13 swapping bodies - nets out
# to the exact same text in B2/B4 as before; it is just Excel re-ordering
# the shared-string table on save, not a value change, so there is nothing
# to replay there.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

# Scroll the view back to the top-left (B1) ...
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1

# ... and leave the selection on D3, matching the saved state in the diff.
$ws.Range("D3").Select() | Out-Null
